$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44162
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("M2").Value = 75
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 18933
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1262
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44162
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 23000
$ws.Range("P3").Value = 23000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 1278

$ws.Range("D4").Value = 44181
$ws.Range("M4").Value = 140
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 17000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 944

$ws.Range("D5").Value = 44176
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("S5").Value = 1111
$ws.Range("T5").Value = 18

$ws.Range("D8").Value = 44160
$ws.Range("M8").Value = 175
$ws.Range("N8").Value = 18000
$ws.Range("P8").Value = 18743
$ws.Range("S8").Value = 1250

$ws.Range("D9").Value = 44167
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("M9").Value = 100
$ws.Range("Q9").Value = '$/caja 15 kilos'
$ws.Range("S9").Value = 1333
$ws.Range("T9").Value = 15

$ws.Range("D10").Value = 44217
$ws.Range("K10").Value = 'Modesto'
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = '$/bandeja 18 kilos'
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44175
$ws.Range("K11").Value = 'Castle Brite'
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("S11").Value = 1111

$ws.Range("D12").Value = 44175
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("S12").Value = 1000

$ws.Range("D13").Value = 44175
$ws.Range("L13").Value = 'Tercera'
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("S13").Value = 778

$ws.Range("D14").Value = 44188
$ws.Range("M14").Value = 35
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("Q14").Value = '$/bandeja 18 kilos'
$ws.Range("S14").Value = 1111
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44174
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 21083
$ws.Range("Q15").Value = '$/bandeja 18 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1171
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44168
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 450
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22444
$ws.Range("Q16").Value = '$/bandeja 18 kilos'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1247
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 44186
$ws.Range("K17").Value = 'Modesto'
$ws.Range("M17").Value = 55
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 20000
$ws.Range("S17").Value = 1111

$ws.Range("D18").Value = 44202
$ws.Range("K18").Value = 'Modesto'
$ws.Range("M18").Value = 25
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 18000
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1200

$ws.Range("D19").Value = 44159
$ws.Range("M19").Value = 85
$ws.Range("Q19").Value = '$/caja 15 kilos'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1333
$ws.Range("T19").Value = 15

$ws.Range("D20").Value = 44166
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 17750
$ws.Range("Q20").Value = '$/caja 15 kilos'
$ws.Range("R20").Value = 'Provincia de Limarí'
$ws.Range("S20").Value = 1183
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 44166
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 15 kilos'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 800
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 44172
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 22000
$ws.Range("Q22").Value = '$/bandeja 18 kilos'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1222
$ws.Range("T22").Value = 18

$ws.Range("D23").Value = 44172
$ws.Range("M23").Value = 65
$ws.Range("Q23").Value = '$/bandeja 18 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'

$ws.Range("D24").Value = 44201
$ws.Range("M24").Value = 45
$ws.Range("N24").Value = 18000
$ws.Range("O24").Value = 18000
$ws.Range("P24").Value = 18000
$ws.Range("Q24").Value = '$/caja 15 kilos'
$ws.Range("S24").Value = 1200
$ws.Range("T24").Value = 15
